# Scenarios.xlsx edit: refine specification of output intervals
#
# - Column G ("SimulationTime") values are changed from a plain number
#   (points per time unit) to a textual triplet specification
#   "<start>, <end>, <resolution>" (multiple intervals separated by ";").
# - A clarifying cell comment (by Pavel Balazki) is added to the header
#   cell G1 explaining the new expected format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# --- Update the SimulationTime values (column G) to the new triplet format ---
$ws.Range("G2").Value = "0, 24, 60"
$ws.Range("G3").Value = "0, 1, 60; 1, 12, 20"
$ws.Range("G4").Value = "0, 12, 20"
$ws.Range("G5").Value = "0, 12, 20"

# --- Add explanatory comment to the header cell G1 ---
$commentText = "Pavel Balazki:" + [char]10 + "Simulation time is defined as time intervals." + [char]10 + "Expected is a triple of values {start, end, resolution}, resolution given in ""points per <time unit>"" as defined in the columne ""SimulationTimeUnit"". Multiple intervals can be separated by a "";"""
$ws.Range("G1").AddComment($commentText) | Out-Null

# --- Restore the final selection on the Scenarios sheet ---
$ws.Range("K22").Select() | Out-Null
